$wb = $excel.ActiveWorkbook

$wsSummary = $wb.Worksheets.Item(1)     # "Итог"
$wsDetail  = $wb.Worksheets.Item(2)     # "Детальный отчет"

# Duplicate the "Итог" sheet twice. Worksheet.Copy keeps all formatting
# (column widths, row heights, cell styles, merged cells, fonts) intact,
# and inserts the copies right after the sheet being copied.
$wsSummary.Copy([Type]::Missing, $wsSummary)
$wsHttp = $wb.Worksheets.Item(2)
$wsHttp.Name = "Итог HTTP"

$wsSummary.Copy([Type]::Missing, $wsHttp)
$wsHttps = $wb.Worksheets.Item(3)
$wsHttps.Name = "Итог HTTPS"

# Match the per-sheet selection / active cell state from the target workbook.
$wsHttp.Range("A1:C2").Select() | Out-Null

$wsHttps.Range("A9").Select() | Out-Null

$wsSummary.Range("A29").Select() | Out-Null
$wsSummary.Activate()

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Index $s.Name
}
